$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2780
$ws1.Range("F4").Value = 1111
$ws1.Range("F5").Value = 20249
$ws1.Range("F7").Value = 2445
$ws1.Range("F11").Value = 715
$ws1.Range("F12").Value = 263
$ws1.Range("F15").Value = 388

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 138

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6032
$ws3.Range("F3").Value = 669
$ws3.Range("F4").Value = 622
$ws3.Range("F5").Value = 1286

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6032
$ws4.Range("F3").Value = 669
$ws4.Range("F4").Value = 622
$ws4.Range("F7").Value = 1286
$ws4.Range("F8").Value = 2780
$ws4.Range("F9").Value = 1111
$ws4.Range("F10").Value = 20249
$ws4.Range("F16").Value = 2445
$ws4.Range("F18").Value = 138
$ws4.Range("F21").Value = 715
$ws4.Range("F22").Value = 263
$ws4.Range("F28").Value = 388
